# Weekly update: insert a new price record as row 13 on the active sheet,
# pushing the existing rows 13:39 down to 14:40 (dimension grows to A1:R40).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 13, shifting rows 13-39 down to 14-40.
$ws.Rows(13).Insert()

# Populate the newly-inserted row 13 with the new weekly record.
$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(13, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(13, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(13, 4).Value = 44565
$ws.Cells.Item(13, 5).Value = 15
$ws.Cells.Item(13, 6).Value = 100112027
$ws.Cells.Item(13, 7).Value = "Melón"
$ws.Cells.Item(13, 8).Value = "Tuna"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 40
$ws.Cells.Item(13, 11).Value = 10000
$ws.Cells.Item(13, 12).Value = 11000
$ws.Cells.Item(13, 13).Value = 10500
$ws.Cells.Item(13, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(13, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 16).Value = 583
$ws.Cells.Item(13, 17).Value = 18
$ws.Cells.Item(13, 18).Value = "Hortaliza"
